# GameData.xlsx — implement buyMode:
#  - Dialog: drop the now-unused "What would you like?" follow-up message
#    (row 12) and the stray leftover test rows (17-20); the shop's first
#    line now fires event 3 (empty-enemy-list guard) via F11.
#  - Items: give the new shopkeeper (IndividualID 207) two items for sale.
#  - Individuals: add the new shopkeeper NPC (ID 207).
#  - Events: add event 3 (the enemyActionMode empty-check / end-start cycle).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dialog")
$ws2 = $wb.Worksheets.Item("Items")
$ws3 = $wb.Worksheets.Item("Individuals")
$ws4 = $wb.Worksheets.Item("Events")

# ---------------------------------------------------------------------
# Dialog (sheet1): remove dialog message 1005 ("What would you like?")
# and the orphaned rows 17-20, then point the shop-intro line at the new
# empty-enemy-list event (F11: 0 -> 3). Delete bottom-up so earlier row
# numbers stay valid.
# ---------------------------------------------------------------------
$ws1.Rows.Item(20).Delete()
$ws1.Rows.Item(19).Delete()
$ws1.Rows.Item(18).Delete()
$ws1.Rows.Item(17).Delete()
$ws1.Rows.Item(12).Delete()
$ws1.Cells.Item(11, 6).Value = 3

# ---------------------------------------------------------------------
# Items (sheet2): new rows for the shopkeeper's wares.
# ---------------------------------------------------------------------
$ws2.Cells.Item(8,1).Value = 207
$ws2.Cells.Item(8,2).Value = 3023
$ws2.Cells.Item(8,3).Value = 1102
$ws2.Cells.Item(8,4).Value = 'i'
$ws2.Cells.Item(8,5).Value = 'p'
$ws2.Cells.Item(8,6).Value = 'l'
$ws2.Cells.Item(8,7).Value = 'c'
$ws2.Cells.Item(8,8).Value = 20
$ws2.Cells.Item(8,9).Value = 255
$ws2.Cells.Item(8,10).Value = 0
$ws2.Cells.Item(8,11).Value = 255
$ws2.Cells.Item(8,12).Value = 'Mana Potion'
$ws2.Cells.Item(8,13).Value = 1
$ws2.Cells.Item(8,14).Value = 2
$ws2.Cells.Item(8,15).Value = 0
$ws2.Cells.Item(8,16).Value = 0
$ws2.Cells.Item(8,17).Value = 0
$ws2.Cells.Item(8,18).Value = 2
$ws2.Cells.Item(8,19).Value = 0
$ws2.Cells.Item(8,20).Value = 0
$ws2.Cells.Item(8,21).Value = 0
$ws2.Cells.Item(8,22).Value = 0
$ws2.Cells.Item(8,23).Value = 0
$ws2.Cells.Item(8,24).Value = 0
$ws2.Cells.Item(8,25).Value = 0
$ws2.Cells.Item(8,26).Value = 0
$ws2.Cells.Item(8,27).Value = 0
$ws2.Cells.Item(8,28).Value = 0
$ws2.Cells.Item(8,29).Value = 0
$ws2.Cells.Item(8,30).Value = 0
$ws2.Cells.Item(8,31).Value = 0
$ws2.Cells.Item(8,32).Value = 0
$ws2.Cells.Item(8,33).Value = 0
$ws2.Cells.Item(8,34).Value = 0
$ws2.Cells.Item(8,35).Value = 0
$ws2.Cells.Item(8,36).Value = 0
$ws2.Cells.Item(8,37).Value = 0
$ws2.Cells.Item(8,38).Value = 0
$ws2.Cells.Item(8,39).Value = 0
$ws2.Cells.Item(8,40).Value = 0
$ws2.Cells.Item(8,41).Value = 'A potion which restores mana.&&(Restores 4 Mana)'

$ws2.Cells.Item(9,1).Value = 207
$ws2.Cells.Item(9,2).Value = 3022
$ws2.Cells.Item(9,3).Value = 1104
$ws2.Cells.Item(9,4).Value = 'i'
$ws2.Cells.Item(9,5).Value = 'p'
$ws2.Cells.Item(9,6).Value = 'l'
$ws2.Cells.Item(9,7).Value = 'c'
$ws2.Cells.Item(9,8).Value = 15
$ws2.Cells.Item(9,9).Value = 255
$ws2.Cells.Item(9,10).Value = 0
$ws2.Cells.Item(9,11).Value = 255
$ws2.Cells.Item(9,12).Value = 'Health Potion'
$ws2.Cells.Item(9,13).Value = 2
$ws2.Cells.Item(9,14).Value = 3
$ws2.Cells.Item(9,15).Value = 0
$ws2.Cells.Item(9,16).Value = 10
$ws2.Cells.Item(9,17).Value = 0
$ws2.Cells.Item(9,18).Value = 0
$ws2.Cells.Item(9,19).Value = 0
$ws2.Cells.Item(9,20).Value = 0
$ws2.Cells.Item(9,21).Value = 0
$ws2.Cells.Item(9,22).Value = 0
$ws2.Cells.Item(9,23).Value = 0
$ws2.Cells.Item(9,24).Value = 0
$ws2.Cells.Item(9,25).Value = 0
$ws2.Cells.Item(9,26).Value = 0
$ws2.Cells.Item(9,27).Value = 0
$ws2.Cells.Item(9,28).Value = 0
$ws2.Cells.Item(9,29).Value = 0
$ws2.Cells.Item(9,30).Value = 0
$ws2.Cells.Item(9,31).Value = 0
$ws2.Cells.Item(9,32).Value = 0
$ws2.Cells.Item(9,33).Value = 0
$ws2.Cells.Item(9,34).Value = 0
$ws2.Cells.Item(9,35).Value = 0
$ws2.Cells.Item(9,36).Value = 0
$ws2.Cells.Item(9,37).Value = 0
$ws2.Cells.Item(9,38).Value = 0
$ws2.Cells.Item(9,39).Value = 0
$ws2.Cells.Item(9,40).Value = 0
$ws2.Cells.Item(9,41).Value = 'A potion which restores health.&&(Restores 10 HP)'

# ---------------------------------------------------------------------
# Individuals (sheet3): new shopkeeper NPC, ID 207.
# ---------------------------------------------------------------------
$ws3.Cells.Item(9,1).Value = 2013
$ws3.Cells.Item(9,2).Value = 207
$ws3.Cells.Item(9,3).Value = 255
$ws3.Cells.Item(9,4).Value = 70
$ws3.Cells.Item(9,5).Value = 255
$ws3.Cells.Item(9,6).Value = 'red_robe'
$ws3.Cells.Item(9,7).Value = 0
$ws3.Cells.Item(9,8).Value = 1
$ws3.Cells.Item(9,9).Value = 2
$ws3.Cells.Item(9,10).Value = 15
$ws3.Cells.Item(9,11).Value = 2
$ws3.Cells.Item(9,12).Value = 10
$ws3.Cells.Item(9,13).Value = 10
$ws3.Cells.Item(9,14).Value = 0
$ws3.Cells.Item(9,15).Value = 5
$ws3.Cells.Item(9,16).Value = 1
$ws3.Cells.Item(9,17).Value = 'MAX'
$ws3.Cells.Item(9,18).Value = 1
$ws3.Cells.Item(9,19).Value = 3
$ws3.Cells.Item(9,20).Value = 0
$ws3.Cells.Item(9,21).Value = 0
$ws3.Cells.Item(9,22).Value = 0
$ws3.Cells.Item(9,23).Value = 0
$ws3.Cells.Item(9,24).Value = 0
$ws3.Cells.Item(9,25).Value = 0
$ws3.Cells.Item(9,26).Value = 0
$ws3.Cells.Item(9,27).Value = 0
$ws3.Cells.Item(9,28).Value = 0
$ws3.Cells.Item(9,29).Value = 0
$ws3.Cells.Item(9,30).Value = 0
$ws3.Cells.Item(9,31).Value = 0
$ws3.Cells.Item(9,32).Value = 1004
$ws3.Cells.Item(9,33).Value = 100

# ---------------------------------------------------------------------
# Events (sheet4): new event 3 - enemyActionMode empty-check guard.
# ---------------------------------------------------------------------
$ws4.Cells.Item(4,1).Value = 3
$ws4.Cells.Item(4,2).Value = 3
$ws4.Cells.Item(4,3).Value = 0
$ws4.Cells.Item(4,4).Value = 207
$ws4.Cells.Item(4,5).Value = 0
$ws4.Cells.Item(4,6).Value = 0
$ws4.Cells.Item(4,7).Value = 0
$ws4.Cells.Item(4,8).Value = 0

# ---------------------------------------------------------------------
# View state: selection per sheet, and Events becomes the active tab.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A11:F11").Select()

$ws2.Activate()
$ws2.Range("A2:AO9").Select()

$ws3.Activate()
$ws3.Range("A9:AG9").Select()

$ws4.Activate()
$ws4.Range("E10").Select()
